$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix formatting of "purpose" column (E) values from "fullRNASEQ" to "fullRNASeq"
$ws.Range("E2:E13").Value = "fullRNASeq"
